# Clean up of driver implementation:
# Insert a new "Type" / "Result" column between the TestCase Name column
# and the old Param1/Values columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B, shifting the former B (Param1) and C (Values)
# columns one place to the right (-> C and D).
$ws.Columns("B").Insert()

# New header + value for the inserted column.
$ws.Range("B1").Value = "Type"
$ws.Range("B2").Value = "Result"

# Match the column width used for the new column (stored width "15").
$ws.Columns("B").ColumnWidth = 14.1667

# Leave the selection on the newly added cell.
$ws.Range("B2").Select()
